# Apply cryptocurrency price/volume updates from the latest GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($range, $value)
    # Force text number-format so numeric-looking strings (e.g. "167.98")
    # are not silently coerced into floating point numbers, then restore
    # the default "Normal" style so no stray formatting is introduced.
    $cell = $ws.Range($range)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-CellText "D2" "67.156.72"
Set-CellText "E2" "  -0.69%  "
Set-CellText "D3" "2.474.39"
Set-CellText "E3" "  -0.72%  "
Set-CellText "D5" "582.13"
Set-CellText "E5" "  -1.47%  "
Set-CellText "D6" "167.98"
Set-CellText "E6" "  -3.42%  "
Set-CellText "E7" "  +0.05%  "
Set-CellText "E8" "  -1.76%  "
Set-CellText "D9" "2.474.72"
Set-CellText "E9" "  -0.71%  "
Set-CellText "E10" "  -3.55%  "
Set-CellText "D11" "0.165"
Set-CellText "E11" "  -0.99%  "
Set-CellText "E12" "  -2.65%  "
Set-CellText "E13" "  -2.50%  "
Set-CellText "B14" "WrappedliquidstakedEther2.0"
Set-CellText "C14" "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-CellText "D14" "2.924.86"
Set-CellText "E14" "  -0.96%  "
Set-CellText "B15" "Avalanche"
Set-CellText "C15" "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-CellText "D15" "25.48"
Set-CellText "E15" "  -3.17%  "
Set-CellText "D16" "66.994.43"
Set-CellText "E16" "  -0.93%  "
Set-CellText "D17" "0.0000169"
Set-CellText "E17" "  -4.07%  "
Set-CellText "D18" "2.487.92"
Set-CellText "E18" "  -0.10%  "
Set-CellText "D19" "11.28"
Set-CellText "E19" "  -4.63%  "
Set-CellText "D20" "7.62"
Set-CellText "E20" "  -4.81%  "
Set-CellText "D21" "357.01"
Set-CellText "E21" "  -2.82%  "
Set-CellText "D22" "4.05"
Set-CellText "E22" "  -1.80%  "
Set-CellText "E23" "  -0.04%  "
Set-CellText "D24" "69.43"
Set-CellText "E25" "  -6.98%  "
Set-CellText "D26" "1.78"
Set-CellText "E26" "  -7.14%  "
Set-CellText "D27" "9.12"
Set-CellText "E27" "  -8.59%  "
Set-CellText "E28" "  +0.18%  "
Set-CellText "E29" "  -0.93%  "
Set-CellText "D30" "0.0₃0905"
Set-CellText "E30" "  -5.56%  "
Set-CellText "D31" "510.38"
Set-CellText "E31" "  -4.07%  "
Set-CellText "D32" "7.79"
Set-CellText "E32" "  -6.41%  "
Set-CellText "E33" "  -4.41%  "
Set-CellText "E34" "  -5.67%  "
Set-CellText "E35" "  +0.04%  "
Set-CellText "E36" "  -6.40%  "
Set-CellText "D37" "158.21"
Set-CellText "E37" "  +0.29%  "
Set-CellText "B38" "WhiteBITCoin"
Set-CellText "C38" "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-CellText "D38" "18.60"
Set-CellText "E38" "  -0.25%  "
Set-CellText "B39" "EthereumClassic"
Set-CellText "C39" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-CellText "D39" "18.45"
Set-CellText "E39" "  -1.33%  "
Set-CellText "E40" "  -5.49%  "
Set-CellText "E42" "  -6.07%  "
Set-CellText "D43" "4.80"
Set-CellText "E43" "  -6.13%  "
Set-CellText "E44" "  -6.43%  "
Set-CellText "E45" "  -6.41%  "
Set-CellText "E46" "  -2.41%  "
Set-CellText "D47" "141.47"
Set-CellText "E47" "  -2.32%  "
Set-CellText "D48" "3.47"
Set-CellText "E48" "  -5.60%  "
Set-CellText "D49" "0.515"
Set-CellText "E49" "  -5.76%  "
Set-CellText "E50" "  -5.44%  "
Set-CellText "D51" "0.0₆0251"
Set-CellText "E51" "  -8.30%  "
